$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-29 Sunday" "2024-09-30 Monday"

Replace-Text "450×6=" "595×9="
Replace-Text "647×7=" "152×6="
Replace-Text "532×4=" "806×6="
Replace-Text "644×9=" "542×8="
Replace-Text "225×6=" "132×5="
Replace-Text "768×4=" "294×4="
Replace-Text "454×4=" "857×4="
Replace-Text "282×6=" "161×5="
Replace-Text "420×8=" "920×8="
Replace-Text "633×2=" "536×5="
Replace-Text "873×8=" "422×8="
Replace-Text "494×9=" "548×6="
Replace-Text "694×8=" "525×5="
Replace-Text "152×2=" "635×5="
Replace-Text "889×2=" "493×5="
Replace-Text "721×5=" "469×2="
Replace-Text "994×9=" "185×3="
Replace-Text "338×5=" "187×8="
Replace-Text "681×6=" "432×9="
Replace-Text "810×8=" "400×7="
Replace-Text "499×2=" "762×9="
Replace-Text "288×6=" "564×8="
Replace-Text "156×9=" "591×3="
Replace-Text "276×4=" "459×3="
Replace-Text "854×2=" "779×3="
